$wb = $excel.ActiveWorkbook

# --- ALC sheet: row 32 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1469
$ws.Range("I32").Value = 1469
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1469
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1143
$ws.Range("N32").ClearContents()

# --- BSM sheet: row 20 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1759.5555
$ws.Range("I20").Value = 1406.7273
$ws.Range("J20").Value = 2314
$ws.Range("K20").Value = 1406.7273
$ws.Range("L20").Value = 2314
$ws.Range("M20").Value = -1159.7273
$ws.Range("N20").Value = -2808

# --- CUL sheet: row 107 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 380.75
$ws.Range("I107").Value = 306.57144
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 919.71432
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 1000.28568
$ws.Range("N107").Value = -6540

# --- GSM sheet: row 12, 21, 30 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H21").Value = 18000
$ws.Range("I21").Value = 18000
$ws.Range("K21").Value = 18000
$ws.Range("M21").Value = -17827

$ws.Range("H30").Value = 18000
$ws.Range("I30").Value = 18000
$ws.Range("K30").Value = 18000
$ws.Range("M30").Value = -17895

# --- LTW sheet: row 46 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5998
$ws.Range("I46").Value = 4016
$ws.Range("K46").Value = 4016
$ws.Range("M46").Value = -3828

# --- WVR sheet: rows 119-141 (fill in new market-price columns) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0

$ws.Range("H122").Value = 4387
$ws.Range("I122").Value = 2420
$ws.Range("J122").Value = 5862.25
$ws.Range("K122").Value = 7260
$ws.Range("L122").Value = 17586.75
$ws.Range("M122").Value = -4810
$ws.Range("N122").Value = -22486.75

$ws.Range("H123").Value = 75000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 75000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800

$ws.Range("H124").Value = 100000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 100000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws.Range("H126").Value = 5549.8
$ws.Range("I126").Value = 3769.6
$ws.Range("J126").Value = 7330
$ws.Range("K126").Value = 11308.8
$ws.Range("L126").Value = 21990
$ws.Range("M126").Value = -8838.8
$ws.Range("N126").Value = -26930

$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 239499.5
$ws.Range("I129").Value = 49999
$ws.Range("J129").Value = 429000
$ws.Range("K129").Value = 49999
$ws.Range("L129").Value = 429000
$ws.Range("M129").Value = -44999
$ws.Range("N129").Value = -439000

$ws.Range("H130").Value = 24000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 24000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 24000
$ws.Range("N130").Value = -34040

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0

$ws.Range("H132").Value = 3483.3333
$ws.Range("I132").Value = 3483.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10449.9999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7919.999899999999

$ws.Range("H133").Value = 120000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 120000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 120000
$ws.Range("N133").Value = -130120

$ws.Range("H135").Value = 39296.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39296.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39296.5
$ws.Range("N135").Value = -49436.5

$ws.Range("H136").Value = 4928.4287
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4928.4287
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 14785.2861
$ws.Range("N136").Value = -19885.2861

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 30000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 30000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

$ws.Range("H141").Value = 299997.5
$ws.Range("I141").Value = 500000
$ws.Range("J141").Value = 99995
$ws.Range("K141").Value = 500000
$ws.Range("L141").Value = 99995
$ws.Range("M141").Value = -494820
$ws.Range("N141").Value = -110355

